$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the "Exp 10" experiment parameters on row 6
$ws.Range("A6").Value = "Exp 10"
$ws.Range("B6").Value = 0.9
$ws.Range("C6").Value = 1
$ws.Range("F6").Value = "Exp 10.png"
